$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column K
$ws.Range("K1").Value = "PRODUCTO"

# Fill K2:K267 with "TRIGO"
$ws.Range("K2:K267").Value = "TRIGO"
